$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$theme = $nm.Theme
$tcs = $theme.ThemeColorScheme
$tcs.Colors(1).RGB = 123456
Write-Host "done"
